$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rename test identifiers: drop "Automobile" from the naming scheme
$ws.Range("A2").Value = "103_TruckInsurance_001_SmokeTest"
$ws.Range("B2").Value = "103_TruckInsurance_001_SmokeTest_FillPageVehicleData"
$ws.Range("C2").Value = "103_TruckInsurance_001_SmokeTest_FillPageInsurantData"
$ws.Range("D2").Value = "103_TruckInsurance_001_SmokeTest_FillPageProductData"
$ws.Range("F2").Value = "103_TruckInsurance_001_SmokeTest_FillPageSendQuote"

# Update the active selection on the sheet
$ws.Range("E17").Select()
